# Auto-generated edit script: update Leve profit calculation values
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 68070.60000000001
$ws.Range("I111").Value = 1266
$ws.Range("J111").Value = 251783.25
$ws.Range("K111").Value = 3798
$ws.Range("L111").Value = 755349.75
$ws.Range("M111").Value = -731
$ws.Range("N111").Value = -761483.75

$ws.Range("H132").Value = 1334.5
$ws.Range("I132").Value = 1398.6923
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 4196.0769
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = -1666.0769
$ws.Range("N132").Value = -6560

$ws.Range("H138").Value = 2457.76
$ws.Range("I138").Value = 836.8823
$ws.Range("J138").Value = 5902.125
$ws.Range("K138").Value = 2510.6469
$ws.Range("L138").Value = 17706.375
$ws.Range("M138").Value = 2629.3531
$ws.Range("N138").Value = -27986.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 122.75
$ws.Range("I4").Value = 130.33333
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 130.33333
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = -14.33332999999999
$ws.Range("N4").Value = -332

$ws.Range("H12").Value = 3519.2
$ws.Range("I12").Value = 2798
$ws.Range("J12").Value = 4000
$ws.Range("K12").Value = 2798
$ws.Range("L12").Value = 4000
$ws.Range("M12").Value = -2625
$ws.Range("N12").Value = -4346

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 540.6667
$ws.Range("I11").Value = 502
$ws.Range("J11").Value = 560
$ws.Range("K11").Value = 502
$ws.Range("L11").Value = 560
$ws.Range("M11").Value = -362
$ws.Range("N11").Value = -840

$ws.Range("H12").Value = 3000
$ws.Range("I12").Value = 2000
$ws.Range("J12").Value = 4000
$ws.Range("K12").Value = 2000
$ws.Range("L12").Value = 4000
$ws.Range("M12").Value = -1832
$ws.Range("N12").Value = -4336

$ws.Range("H107").Value = 942.0833
$ws.Range("I107").Value = 933.8095
$ws.Range("K107").Value = 933.8095
$ws.Range("M107").Value = 986.1905

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 12780.2
$ws.Range("J4").Value = 3633.3333
$ws.Range("L4").Value = 3633.3333
$ws.Range("N4").Value = -3857.3333

$ws.Range("H58").Value = 1590.381
$ws.Range("I58").Value = 1142.7142
$ws.Range("J58").Value = 2485.7144
$ws.Range("K58").Value = 1142.7142
$ws.Range("L58").Value = 2485.7144
$ws.Range("M58").Value = -939.7141999999999
$ws.Range("N58").Value = -2891.7144

$ws.Range("H107").Value = 233.21951
$ws.Range("I107").Value = 169.76923
$ws.Range("K107").Value = 169.76923
$ws.Range("M107").Value = 1750.23077

$ws.Range("H135").Value = 37837.5
$ws.Range("J135").Value = 37837.5
$ws.Range("L135").Value = 37837.5
$ws.Range("N135").Value = -47977.5

$ws.Range("H136").Value = 1590.381
$ws.Range("I136").Value = 1142.7142
$ws.Range("J136").Value = 2485.7144
$ws.Range("K136").Value = 3428.1426
$ws.Range("L136").Value = 7457.1432
$ws.Range("M136").Value = -878.1425999999997
$ws.Range("N136").Value = -12557.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H126").Value = 6913.3335
$ws.Range("I126").Value = 6500
$ws.Range("J126").Value = 7120
$ws.Range("K126").Value = 19500
$ws.Range("L126").Value = 21360
$ws.Range("M126").Value = -14560
$ws.Range("N126").Value = -31240

$ws.Range("H141").Value = 27720.25
$ws.Range("I141").Value = 33960.332
$ws.Range("K141").Value = 101880.996
$ws.Range("M141").Value = -96700.99600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 11001.333
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 11001.333
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 11001.333
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -11339.333

$ws.Range("H13").Value = 10000
$ws.Range("J13").Value = 10000
$ws.Range("L13").Value = 10000
$ws.Range("M13").Value = -10278

$ws.Range("H136").Value = 37933.332
$ws.Range("J136").Value = 37933.332
$ws.Range("L136").Value = 113799.996
$ws.Range("N136").Value = -118899.996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 62501908
$ws.Range("I40").Value = 100001640
$ws.Range("K40").Value = 100001640
$ws.Range("M40").Value = -100001504

$ws.Range("H136").Value = 6196.256
$ws.Range("I136").Value = 5922.8965
$ws.Range("J136").Value = 6762.5
$ws.Range("K136").Value = 17768.6895
$ws.Range("L136").Value = 20287.5
$ws.Range("M136").Value = -15218.6895
$ws.Range("N136").Value = -25387.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 7000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 7000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 7000
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -7280

$ws.Range("H10").Value = 12000
$ws.Range("J10").Value = 12000
$ws.Range("L10").Value = 12000
$ws.Range("N10").Value = -12338

$ws.Range("H13").Value = 3006
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 3006
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 3006
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -3286

$ws.Range("H81").Value = 3001
$ws.Range("I81").Value = 3001
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 6002
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -4941
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 3001
$ws.Range("I84").Value = 3001
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 30010
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -24706
$ws.Range("N84").ClearContents()

$ws.Range("H132").Value = 1190.2322
$ws.Range("I132").Value = 782.3171
$ws.Range("J132").Value = 2305.2
$ws.Range("K132").Value = 2346.9513
$ws.Range("L132").Value = 6915.599999999999
$ws.Range("M132").Value = 183.0487000000003
$ws.Range("N132").Value = -11975.6

Write-Host "Applied 159 cell updates across 27 rows"